$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Steps")
$ws.Activate()

# Add one more test step row (row 10), matching the formatting of the row above it
$ws.Range("A9:E9").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122)

$ws.Range("A10").Value = "searchListingPageTestCases"
$ws.Range("C10").Value = "addProduct"
$ws.Range("D10").Value = "search_list_projectNames|search_list_addBtn"
$ws.Range("E10").Value = "Apple - Royal Gala"

# Update the active selection/view to the new bottom of the sheet
$ws.Range("C13").Select()
